# Clean the dataset: the rows that were previously labelled with the
# "cleaned" class (rows 52-117 of the Constraint_Train sheet) are wiped
# of their content, leaving just the empty, styled cells behind - matching
# the rest of the already-empty tail of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("A52:C117")
$target.ClearContents()

# Reflect the author's on-screen state after the edit: the view is
# scrolled down a bit further and the cleared block stays selected.
$ws.Range("A104").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 104
$aw.ScrollColumn = 1
$ws.Range("A52:C117").Select()
